$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextCell "D2" "42.701.57"
Set-TextCell "E2" "  -0.47%  "
Set-TextCell "D3" "2.300.49"
Set-TextCell "E3" "  -0.38%  "
Set-TextCell "E4" "  +0.13%  "
Set-TextCell "D5" "311.45"
Set-TextCell "E5" "  -2.39%  "
Set-TextCell "D6" "104.85"
Set-TextCell "E6" "  +0.02%  "
Set-TextCell "D7" "0.625"
Set-TextCell "E7" "  -0.99%  "
Set-TextCell "E8" "  +0.01%  "
Set-TextCell "D9" "0.605"
Set-TextCell "E9" "  -0.90%  "
Set-TextCell "D10" "39.75"
Set-TextCell "E10" "  -1.32%  "
Set-TextCell "E11" "  -0.92%  "
Set-TextCell "E12" "  -4.28%  "
Set-TextCell "E13" "  +0.15%  "
Set-TextCell "D14" "0.986"
Set-TextCell "E14" "  +0.67%  "
Set-TextCell "D15" "2.782.18"
Set-TextCell "E15" "  +4.76%  "
Set-TextCell "E16" "  -0.37%  "
Set-TextCell "D17" "2.292.41"
Set-TextCell "E17" "  -1.58%  "
Set-TextCell "D18" "42.824.12"
Set-TextCell "E18" "  -0.01%  "
Set-TextCell "D19" "7.30"
Set-TextCell "E19" "  -3.29%  "
Set-TextCell "D20" "13.62"
Set-TextCell "E20" "  +0.27%  "
Set-TextCell "D21" "0.0000105"
Set-TextCell "E21" "  -1.61%  "
Set-TextCell "D22" "73.39"
Set-TextCell "E22" "  -0.75%  "
Set-TextCell "D23" "3.46"
Set-TextCell "E23" "  -2.97%  "
Set-TextCell "D24" "269.31"
Set-TextCell "E24" "  -0.87%  "
Set-TextCell "D25" "2.22"
Set-TextCell "E25" "  -2.62%  "
Set-TextCell "E26" "  +0.57%  "
Set-TextCell "D27" "10.85"
Set-TextCell "E27" "  -1.02%  "
Set-TextCell "D28" "7.14"
Set-TextCell "E28" "  +13.60%  "
Set-TextCell "E29" "  -1.15%  "
Set-TextCell "D30" "22.32"
Set-TextCell "E30" "  -1.87%  "
Set-TextCell "D31" "36.19"
Set-TextCell "D32" "164.84"
Set-TextCell "E32" "  -0.78%  "
Set-TextCell "D33" "0.0856"
Set-TextCell "E33" "  -4.22%  "
Set-TextCell "B34" "Stellar"
Set-TextCell "C34" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D34" "0.130"
Set-TextCell "E34" "  -2.15%  "
Set-TextCell "B35" "WEMIXToken"
Set-TextCell "C35" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D35" "2.62"
Set-TextCell "E35" "  +2.79%  "
Set-TextCell "E36" "  -3.73%  "
Set-TextCell "E37" "  -1.57%  "
Set-TextCell "E38" "  -2.29%  "
Set-TextCell "E39" "  +2.10%  "
Set-TextCell "D40" "3.62"
Set-TextCell "E40" "  -2.89%  "
Set-TextCell "D41" "110.44"
Set-TextCell "E41" "  +9.30%  "
Set-TextCell "E42" "  -0.27%  "
Set-TextCell "D43" "71.04"
Set-TextCell "E43" "  +0.25%  "
Set-TextCell "D44" "0.227"
Set-TextCell "E44" "  -0.15%  "
Set-TextCell "D45" "1.01"
Set-TextCell "E45" "  +0.27%  "
Set-TextCell "D46" "12.27"
Set-TextCell "E46" "  -0.92%  "
Set-TextCell "D47" "1.731.03"
Set-TextCell "E47" "  +8.24%  "
Set-TextCell "D48" "110.64"
Set-TextCell "E48" "  -3.72%  "
Set-TextCell "D49" "77.30"
Set-TextCell "E49" "  -7.11%  "
Set-TextCell "D50" "8.64"
Set-TextCell "E50" "  -2.93%  "
Set-TextCell "E51" "  -3.48%  "
